$d = $word.ActiveDocument

# --- Hunk 1: title paragraph - drop w:hint="eastAsia" from the paragraph-mark rFonts ---
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("CS 611 Legend_Of_Valor", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "Hunk1 anchor not found" }
$para1 = $rng1.Paragraphs(1)
$r1 = $para1.Range
$r1.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="1"/><w:jc w:val="center"/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">CS 611 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Legend_Of_Valor</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

# --- Hunk 2: "levelUp() method." - drop gramStart/gramEnd proofErr, merge "(" + ") method." runs ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("levelUp", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Hunk2 anchor not found" }
$para2 = $rng2.Paragraphs(1)
$r2 = $para2.Range
$r2.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">The classes for different types of heroes. Each of them has their own </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t>levelUp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t>() method.</w:t></w:r></w:p>')

# --- Hunk 3: "formulas are able to" - drop gramStart/gramEnd proofErr, merge runs ---
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("are able to", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) { throw "Hunk3 anchor not found" }
$para3 = $rng3.Paragraphs(1)
$r3 = $para3.Range
$r3.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">The class of most formulas. We </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>calculate</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> the damage, dodge chance, Exp and Gold after battle </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>……</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> in </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>this</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> class. Most of </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t>formulas are able to</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> accept different multi</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">ply powering to fit in </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>different</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> logics.</w:t></w:r></w:p>')

# --- Hunk 4: replace the first trailing empty paragraph with the new Patterns / design evaluation sections ---
$rng4 = $d.Content
$found4 = $rng4.Find.Execute("4. workFlow", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found4) { throw "Hunk4 anchor not found" }
$workflowPara = $rng4.Paragraphs(1)
$targetPara = $workflowPara.Next().Next()
$r4 = $targetPara.Range
$r4.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="2"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Patterns</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">1. We use </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/><w:b/><w:bCs/></w:rPr><w:t>Factory</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> pattern</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> to handle the heroes loading process. If we want to add a new type of hero, we only need to add the .txt file and a new class for the hero.</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:br/><w:t xml:space="preserve">2. We use </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/><w:b/><w:bCs/></w:rPr><w:t>F</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:b/><w:bCs/></w:rPr><w:t>acade pattern</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>CharacterMaps</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>itemLists</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>. Both of them provide a unified interface to get the data from .txt files.</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:br/><w:t xml:space="preserve">3. We use </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/><w:b/><w:bCs/></w:rPr><w:t>S</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:b/><w:bCs/></w:rPr><w:t>ingleton pattern</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>AudioPlayer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">. We only need one instance of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>AudioPlayer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> to play music.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="2"/></w:pPr><w:r><w:t>design evaluation</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t>1. We created some tools to pick random items from a map or a list. Because we load and put all items/heroes/monsters into maps or lists. In this way, we are able to get a random item at any time.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">2. We create </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t>CharaterMaps</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t>itemLists</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> to handle all characters and items, which is a </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>façade</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> pattern. In this way, w</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>e can manage multiple interfaces more efficiently without disorganized.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">3. </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">All cells will not contain any sign of heroes or monsters. Pieces of heroes and monsters does </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>“</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t>move</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>”</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> from one cell to another cell. We created stacks for every cell to contain pieces. In this way, we can put multiple pieces in one cell.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t>4. We picked parts of battle logics from Heroes and Monsters. To be specific, w</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">e </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>spli</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">t </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>the original turn-based combat logic into single attacks.</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> So, both monsters and heroes can only attack once in their turn.</w:t></w:r></w:p>')

Write-Output "edit complete"
